$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updates - repulled data / mean calculation
$ws.Range("F3").Value = 10
$ws.Range("F4").Value = -7
$ws.Range("F5").Value = -6
$ws.Range("F7").Value = -7
$ws.Range("F8").Value = -7
$ws.Range("F10").Value = -6
$ws.Range("F12").Value = 1
